$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new values look like plain numbers need to be forced to text
# so Excel does not auto-convert them to numeric values (matching original inlineStr text cells).
$textCells = @("D5", "D6", "D9", "D10", "D12", "D13", "D16", "D18", "D19", "D20", "D21", "D25", "D26", "D27", "D28", "D34", "D35", "D36", "D37", "D40", "D42", "D47", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "595.04"
$ws.Range("D6").Value = "150.70"
$ws.Range("D9").Value = "0.114"
$ws.Range("D10").Value = "5.79"
$ws.Range("D12").Value = "0.152"
$ws.Range("D13").Value = "27.93"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("D18").Value = "12.26"
$ws.Range("D19").Value = "4.79"
$ws.Range("D20").Value = "348.33"
$ws.Range("D21").Value = "6.99"
$ws.Range("D25").Value = "9.26"
$ws.Range("D26").Value = "1.68"
$ws.Range("D27").Value = "8.49"
$ws.Range("D28").Value = "549.48"
$ws.Range("D34").Value = "5.54"
$ws.Range("D35").Value = "6.15"
$ws.Range("D36").Value = "164.40"
$ws.Range("D37").Value = "0.418"
$ws.Range("D40").Value = "19.69"
$ws.Range("D42").Value = "166.89"
$ws.Range("D47").Value = "0.634"
$ws.Range("D50").Value = "19.33"

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining text-like cells (percentages, and D-values with multiple dots/digit groups)
$ws.Range("D2").Value = "63.741.52"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.616.55"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("E10").Value = "  +3.37%  "
$ws.Range("E11").Value = "  +2.84%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("D14").Value = "3.087.26"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "63.581.30"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("E16").Value = "  +10.45%  "
$ws.Range("D17").Value = "2.632.38"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("E24").Value = "  -2.24%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +4.61%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").Value = "0.0₃0887"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("E34").Value = "  +5.86%  "
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  +4.52%  "
$ws.Range("E44").Value = "  +9.79%  "
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("E46").Value = "  +8.03%  "
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  +17.65%  "
